$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-cell format (style s="3": date number format + center/wrap)
# from the last existing row (A28) so new date cells reuse the same style
# instead of Excel fabricating a brand-new cellXf.

# Row 29
$ws.Range("A28").Copy() | Out-Null
$ws.Range("A29").PasteSpecial(-4122) | Out-Null
$ws.Range("A29").Value = 41252
$ws.Range("C29").Value = "Bacon cheese and green pepper omelette; Bagel with butter; Water"
$ws.Range("G29").Value = "Five Guys Burger;Fries;Water"
$ws.Range("I29").Value = "Trail Mix"
$ws.Rows("29").RowHeight = 90

# Row 30
$ws.Range("A28").Copy() | Out-Null
$ws.Range("A30").PasteSpecial(-4122) | Out-Null
$ws.Range("A30").Value = 41253
$ws.Range("G30").Value = "Veggie burger; Sweet potatoe fries"
$ws.Range("I30").Value = "Chocolate"
$ws.Rows("30").RowHeight = 45

# Row 31
$ws.Range("A28").Copy() | Out-Null
$ws.Range("A31").PasteSpecial(-4122) | Out-Null
$ws.Range("A31").Value = 41254
$ws.Range("C31").Value = "Bacon cheese and green pepper omelette; Bagel with butter; Water"
$ws.Range("G31").Value = "Turckey and cheese sandwich"
$ws.Range("H31").Value = "Oreos; Dorritos"
$ws.Rows("31").RowHeight = 90

# Row 32
$ws.Range("A28").Copy() | Out-Null
$ws.Range("A32").PasteSpecial(-4122) | Out-Null
$ws.Range("A32").Value = 41255
$ws.Range("B32").Value = "Scrambled eggs; Bacon; Bagel with butter; Water"
$ws.Range("C32").Value = " "
$ws.Range("E32").Value = "Pulled pork wrap; Water"
$ws.Range("G32").Value = "Chipoltle Burrito; Chips and Guac;Water"
$ws.Rows("32").RowHeight = 60

# Row 33
$ws.Range("A28").Copy() | Out-Null
$ws.Range("A33").PasteSpecial(-4122) | Out-Null
$ws.Range("A33").Value = 41256
$ws.Range("B33").Value = "Scrambled eggs; Sausage; Bagel with butter; Waffle; Butter"
$ws.Range("F33").Value = "Trail mix"
$ws.Range("G33").Value = "Pasta with peas, corn, carrots, chicken, and alfredo sauce; Water"
$ws.Range("H33").Value = "Rice Crispie Treat"
$ws.Range("I33").Value = "Almonds"
$ws.Rows("33").RowHeight = 75

# Row 34
$ws.Range("A28").Copy() | Out-Null
$ws.Range("A34").PasteSpecial(-4122) | Out-Null
$ws.Range("A34").Value = 41257
$ws.Range("C34").Value = "Scrambled eggs; Bacon; French toast; Water"
$ws.Range("G34").Value = "Pulled chicken; Rice; Naan; Peas and mushrooms"
$ws.Range("H34").Value = "Lamb over rice"
$ws.Rows("34").RowHeight = 45

# Row 35
$ws.Range("A28").Copy() | Out-Null
$ws.Range("A35").PasteSpecial(-4122) | Out-Null
$ws.Range("A35").Value = 41258
$ws.Range("C35").Value = "Bagel with cream cheese"
$ws.Range("G35").Value = "Chicken with rice and veggies"
$ws.Range("H35").Value = "Gummy Bears"
$ws.Range("I35").Value = "Trail mix"
$ws.Rows("35").RowHeight = 30

# Row 36
$ws.Range("A28").Copy() | Out-Null
$ws.Range("A36").PasteSpecial(-4122) | Out-Null
$ws.Range("A36").Value = 41259
$ws.Range("G36").Value = "Pasta with chicken, veggies, and alfredo sauce"
$ws.Rows("36").RowHeight = 60

# Row 37
$ws.Range("A28").Copy() | Out-Null
$ws.Range("A37").PasteSpecial(-4122) | Out-Null
$ws.Range("A37").Value = 41260
$ws.Range("B37").Value = "Scrambled eggs; Bagel with butter; French toast; Water"
$ws.Range("E37").Value = "Pasta with peas, carrots, green & red peppers, broccoli, chicken, and alfredo sauce; Cookies; Water"
$ws.Range("G37").Value = "Pork; Mashed potatoes; Green beans; Pasta; Water"
$ws.Range("I37").Value = "Cookies"
$ws.Rows("37").RowHeight = 105

# Row 38
$ws.Range("A28").Copy() | Out-Null
$ws.Range("A38").PasteSpecial(-4122) | Out-Null
$ws.Range("A38").Value = 41261
$ws.Range("E38").Value = "Meatball sub; Salad; Chocolate cake; Water"
$ws.Range("G38").Value = "Chicken Parm; Water"
$ws.Rows("38").RowHeight = 45

# Row 39
$ws.Range("A28").Copy() | Out-Null
$ws.Range("A39").PasteSpecial(-4122) | Out-Null
$ws.Range("A39").Value = 41262
$ws.Range("C39").Value = "Scrambled eggs; Bacon; Oatmeal; Water"
$ws.Range("G39").Value = "Baked Chicken;Mashed potatoes Veggies; Water"
$ws.Range("I39").Value = "Chocolate croissant"
$ws.Rows("39").RowHeight = 60

# Update the view: scrolled down, zoomed to 85%, with A39 selected
$excel.ActiveWindow.Zoom = 85
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A39").Select() | Out-Null
